$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet so the shared-string table / dimension / rows
# rebuild themselves from scratch for the new layout.
$ws.Cells.Clear() | Out-Null

# ---------------------------------------------------------------------
# Top header block (rows 1-7)
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Item Name (Michael Edgar)"
$ws.Range("B1").Value = "Game Manager"

$ws.Range("A2").Value = "Describe Role in game"

$ws.Range("A5").Value = "Internal Functionality"
$ws.Range("B5").Value = "Text Description"

$ws.Range("B2").Value = "The role of the game manager is to keep track of score values, death parmeters, level diffculty."

$ws.Range("B3").Value = "It also instansiates the game objects and asks the world where to spawn them."

$ws.Range("A6").Value = "Spawn Items"
$ws.Range("B6").Value = "Spawn the eggs, ice blocks, enemies, rocks, player, score and pickup items at the start of the level"

$ws.Range("A7").Value = "Set position"
$ws.Range("B7").Value = "Set the position of the items by communicating with the World"

# ---------------------------------------------------------------------
# "External Outgoing" table (rows 10-16)
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "External Outgoing"
$ws.Range("B10").Value = "Text Description"
$ws.Range("C10").Value = "Parameters"
$ws.Range("D10").Value = "Communication with?"

$ws.Range("A11").Value = "Value Display"
$ws.Range("B11").Value = "Communicate with the pop up score to tell it what value to display when  the player and enemy collide."

$ws.Range("A12").Value = "Player Death"

$ws.Range("A13").Value = "Difficulty"
$ws.Range("B13").Value = "Communicate with the enemy to increase movement speed/difficulty."

$ws.Range("A15").Value = "Egg spawn"

$excel.Union($ws.Range("D15"), $ws.Range("D16")).Value = "Egg"

$ws.Range("D11").Value = "Pop Up Score"

$ws.Range("B12").Value = "Communicate with the player when it collides with an enemy to despawn."

$excel.Union($ws.Range("D12"), $ws.Range("D25")).Value = "Player"

$ws.Range("A14").Value = "Enemy Death"
$ws.Range("B14").Value = "Communicate with the enemy when it collides with an ice block to tell it to despawn."

$ws.Range("B15").Value = "Communicate with the egg to tell it when to hatch/spawn an enemy."

$excel.Union($ws.Range("D13"), $ws.Range("D14")).Value = "Enemy"

$ws.Range("A16").Value = "Egg Death"
$ws.Range("B16").Value = "Communicate with the egg when something collides with it to tell it to despawn."

# ---------------------------------------------------------------------
# "External Incoming" table (rows 23-27)
# ---------------------------------------------------------------------
$ws.Range("A23").Value = "External Incoming"
$ws.Range("B23").Value = "Text Description"
$ws.Range("C23").Value = "Return "
$ws.Range("D23").Value = "Communication with?"

$ws.Range("A24").Value = "Score update"
$ws.Range("D24").Value = "Pop up score"
$ws.Range("B24").Value = "Receive score update from the player when the player interacts with an item."

$ws.Range("A25").Value = "Win/Lose Screen"
$ws.Range("B25").Value = "Receive update when player dies to display lose screen or when the player meets the win condition"

$ws.Range("B26").Value = "to display the win screen."

$ws.Range("A27").Value = "Name of InterFace "
$ws.Range("B27").Value = "https://unity3d.com/learn/tutorials/topics/scripting/interfaces"

# Apply the built-in "Hyperlink" look to B27 (underline + theme colour) without
# leaving a real hyperlink relationship behind.
$hlink = $ws.Hyperlinks.Add($ws.Range("B27"), "https://unity3d.com/learn/tutorials/topics/scripting/interfaces")
$ws.Hyperlinks.Delete() | Out-Null

# ---------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 26.140625
$ws.Columns.Item(2).ColumnWidth = 91.42578125
$ws.Columns.Item(3).ColumnWidth = 66.42578125
$ws.Columns.Item(4).ColumnWidth = 107.85546875

# ---------------------------------------------------------------------
# View / selection
# ---------------------------------------------------------------------
$ws.Range("C24").Select() | Out-Null

# ---------------------------------------------------------------------
# Page setup
# ---------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
